$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5843.8887
$ws.Range("I18").Value = 6799.5
$ws.Range("J18").Value = 3932.6667
$ws.Range("K18").Value = 6799.5
$ws.Range("L18").Value = 3932.6667
$ws.Range("M18").Value = -6515.5
$ws.Range("N18").Value = -4500.6667
$ws.Range("H32").Value = 3173.8
$ws.Range("I32").Value = 2232.8333
$ws.Range("J32").Value = 4585.25
$ws.Range("K32").Value = 2232.8333
$ws.Range("L32").Value = 4585.25
$ws.Range("M32").Value = -1906.8333
$ws.Range("N32").Value = -5237.25
$ws.Range("H103").Value = 790.6111
$ws.Range("I103").Value = 478.8
$ws.Range("J103").Value = 910.53845
$ws.Range("K103").Value = 1436.4
$ws.Range("L103").Value = 2731.61535
$ws.Range("M103").Value = -850.4000000000001
$ws.Range("N103").Value = -3903.61535
$ws.Range("H107").Value = 36540236
$ws.Range("I107").Value = 17858626
$ws.Range("J107").Value = 58335450
$ws.Range("K107").Value = 17858626
$ws.Range("L107").Value = 58335450
$ws.Range("M107").Value = -17856706
$ws.Range("N107").Value = -58339290
$ws.Range("H135").Value = 714876.9
$ws.Range("I135").Value = 714876.9
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6433892.100000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6431357.100000001
$ws.Range("H138").Value = 6110.5293
$ws.Range("I138").Value = 1827.5883
$ws.Range("J138").Value = 10393.471
$ws.Range("K138").Value = 5482.7649
$ws.Range("L138").Value = 31180.413
$ws.Range("M138").Value = -342.7649000000001
$ws.Range("N138").Value = -41460.413
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -820
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4188.051
$ws.Range("I61").Value = 2310.6047
$ws.Range("J61").Value = 9233.6875
$ws.Range("K61").Value = 2310.6047
$ws.Range("L61").Value = 9233.6875
$ws.Range("M61").Value = -2098.6047
$ws.Range("N61").Value = -9657.6875
$ws.Range("H122").Value = 14631.895
$ws.Range("I122").Value = 19773.363
$ws.Range("J122").Value = 7562.375
$ws.Range("K122").Value = 59320.08900000001
$ws.Range("L122").Value = 22687.125
$ws.Range("M122").Value = -56870.08900000001
$ws.Range("N122").Value = -27587.125
$ws.Range("H136").Value = 4188.051
$ws.Range("I136").Value = 2310.6047
$ws.Range("J136").Value = 9233.6875
$ws.Range("K136").Value = 6931.8141
$ws.Range("L136").Value = 27701.0625
$ws.Range("M136").Value = -4381.8141
$ws.Range("N136").Value = -32801.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 28433962
$ws.Range("I86").Value = 8965804
$ws.Range("J86").Value = 62503240
$ws.Range("K86").Value = 8965804
$ws.Range("L86").Value = 62503240
$ws.Range("M86").Value = -8964681
$ws.Range("N86").Value = -62505486
$ws.Range("H89").Value = 28433962
$ws.Range("I89").Value = 8965804
$ws.Range("J89").Value = 62503240
$ws.Range("K89").Value = 44829020
$ws.Range("L89").Value = 312516200
$ws.Range("M89").Value = -44823404
$ws.Range("N89").Value = -312527432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4887
$ws.Range("H6").Value = 666
$ws.Range("I6").Value = 666
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 666
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -553
$ws.Range("H50").Value = 49999
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 49999
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 49999
$ws.Range("M50").Value = ""
$ws.Range("N50").Value = -51249
$ws.Range("H51").Value = 51249.75
$ws.Range("I51").Value = 49999
$ws.Range("J51").Value = 51666.668
$ws.Range("K51").Value = 49999
$ws.Range("L51").Value = 51666.668
$ws.Range("M51").Value = -49263
$ws.Range("N51").Value = -53138.668
$ws.Range("H58").Value = 10644755
$ws.Range("I58").Value = 25003282
$ws.Range("J58").Value = 8807.777
$ws.Range("K58").Value = 25003282
$ws.Range("L58").Value = 8807.777
$ws.Range("M58").Value = -25003079
$ws.Range("N58").Value = -9213.777
$ws.Range("H60").Value = 21846.334
$ws.Range("I60").Value = 6062
$ws.Range("J60").Value = 37630.668
$ws.Range("K60").Value = 6062
$ws.Range("L60").Value = 37630.668
$ws.Range("M60").Value = -5551
$ws.Range("N60").Value = -38652.668
$ws.Range("H61").Value = 51249.75
$ws.Range("I61").Value = 49999
$ws.Range("J61").Value = 51666.668
$ws.Range("K61").Value = 49999
$ws.Range("L61").Value = 51666.668
$ws.Range("M61").Value = -49651
$ws.Range("N61").Value = -52362.668
$ws.Range("H136").Value = 10644755
$ws.Range("I136").Value = 25003282
$ws.Range("J136").Value = 8807.777
$ws.Range("K136").Value = 75009846
$ws.Range("L136").Value = 26423.331
$ws.Range("M136").Value = -75007296
$ws.Range("N136").Value = -31523.331
$ws.Range("H141").Value = 61427.332
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 61427.332
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 61427.332
$ws.Range("N141").Value = -71787.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 362.5
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -1424
$ws.Range("H92").Value = 10990696
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 10990696
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 32972088
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = -32974584
$ws.Range("H141").Value = 8716.5
$ws.Range("I141").Value = 2566.5
$ws.Range("J141").Value = 16916.5
$ws.Range("K141").Value = 7699.5
$ws.Range("L141").Value = 50749.5
$ws.Range("M141").Value = -2519.5
$ws.Range("N141").Value = -61109.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3960.1667
$ws.Range("I16").Value = 3062.182
$ws.Range("J16").Value = 5371.2856
$ws.Range("K16").Value = 3062.182
$ws.Range("L16").Value = 5371.2856
$ws.Range("M16").Value = -2892.182
$ws.Range("N16").Value = -5711.2856
$ws.Range("H22").Value = 2510.75
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 2903.2222
$ws.Range("K22").Value = 1333.3334
$ws.Range("L22").Value = 2903.2222
$ws.Range("M22").Value = -1038.3334
$ws.Range("N22").Value = -3493.2222
$ws.Range("H27").Value = 2510.75
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 2903.2222
$ws.Range("K27").Value = 1333.3334
$ws.Range("L27").Value = 2903.2222
$ws.Range("M27").Value = -1226.3334
$ws.Range("N27").Value = -3117.2222
$ws.Range("H87").Value = 70000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 70000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 70000
$ws.Range("N87").Value = -72246
$ws.Range("H88").Value = 39885
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 39885
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 39885
$ws.Range("N88").Value = -40741
$ws.Range("H90").Value = 70000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 70000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 210000
$ws.Range("N90").Value = -221232
$ws.Range("H91").Value = 39885
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 39885
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 39885
$ws.Range("N91").Value = -42849
$ws.Range("H136").Value = 11497.954
$ws.Range("I136").Value = 3488.75
$ws.Range("J136").Value = 13277.777
$ws.Range("K136").Value = 10466.25
$ws.Range("L136").Value = 39833.331
$ws.Range("M136").Value = -7916.25
$ws.Range("N136").Value = -44933.331
$ws.Range("H138").Value = 75000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 75000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1064.7142
$ws.Range("I107").Value = 969.1111
$ws.Range("J107").Value = 1236.8
$ws.Range("K107").Value = 2907.3333
$ws.Range("L107").Value = 3710.4
$ws.Range("M107").Value = -987.3332999999998
$ws.Range("N107").Value = -7550.4
